$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.835941000000001
$ws.Range("H2").Value = 17.507823
$ws.Range("I2").Value = 0.03643643319117328
$ws.Range("J2").Value = 0.03643643319117327
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 15.89711301056834
$ws.Range("R2").Value = 143.074017095115
$ws.Range("S2").Value = 0.001689011033371737
$ws.Range("T2").Value = 0.001689011033371736

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.835941000000001
$ws.Range("H3").Value = 17.507823
$ws.Range("I3").Value = 0.03643643319117328
$ws.Range("J3").Value = 0.03643643319117327
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 237.0596039534574
$ws.Range("R3").Value = 2133.536435581117
$ws.Range("S3").Value = 0.02518672958907331
$ws.Range("T3").Value = 0.0251867295890733

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.835941000000001
$ws.Range("H4").Value = 17.507823
$ws.Range("I4").Value = 0.03643643319117328
$ws.Range("J4").Value = 0.03643643319117327
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 89.98603752218501
$ws.Range("R4").Value = 809.874337699665
$ws.Range("S4").Value = 0.009560692568728232
$ws.Range("T4").Value = 0.00956069256872823

# Row 5
$ws.Range("G5").Value = 17.50798033333334
$ws.Range("H5").Value = 52.52394100000001
$ws.Range("I5").Value = 0.1093102818770573
$ws.Range("J5").Value = 0.1093102818770573
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 47.69176760796723
$ws.Range("R5").Value = 429.2259084717051
$ws.Range("S5").Value = 0.005067078634800348
$ws.Range("T5").Value = 0.005067078634800347

# Row 6
$ws.Range("G6").Value = 17.50798033333334
$ws.Range("H6").Value = 52.52394100000001
$ws.Range("I6").Value = 0.1093102818770573
$ws.Range("J6").Value = 0.1093102818770573
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("Q6").Value = 711.1852028395971
$ws.Range("R6").Value = 6400.666825556374
$ws.Range("S6").Value = 0.0755608677857573
$ws.Range("T6").Value = 0.07556086778575727

# Row 7
$ws.Range("G7").Value = 17.50798033333334
$ws.Range("H7").Value = 52.52394100000001
$ws.Range("I7").Value = 0.1093102818770573
$ws.Range("J7").Value = 0.1093102818770573
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 269.9605385340617
$ws.Range("R7").Value = 2429.644846806556
$ws.Range("S7").Value = 0.02868233545649965
$ws.Range("T7").Value = 0.02868233545649965

# Row 8
$ws.Range("G8").Value = 136.8238143333333
$ws.Range("H8").Value = 410.471443
$ws.Range("I8").Value = 0.8542532849317694
$ws.Range("J8").Value = 0.8542532849317694
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 372.7082982836906
$ws.Range("R8").Value = 3354.374684553216
$ws.Range("S8").Value = 0.03959891507419384
$ws.Range("T8").Value = 0.03959891507419384

# Row 9
$ws.Range("G9").Value = 136.8238143333333
$ws.Range("H9").Value = 410.471443
$ws.Range("I9").Value = 0.8542532849317694
$ws.Range("J9").Value = 0.8542532849317694
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 5557.869628438907
$ws.Range("R9").Value = 50020.82665595016
$ws.Range("S9").Value = 0.5905036416508046
$ws.Range("T9").Value = 0.5905036416508045

# Row 10
$ws.Range("G10").Value = 136.8238143333333
$ws.Range("H10").Value = 410.471443
$ws.Range("I10").Value = 0.8542532849317694
$ws.Range("J10").Value = 0.8542532849317694
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 2109.725387992752
$ws.Range("R10").Value = 18987.52849193477
$ws.Range("S10").Value = 0.224150728206771
$ws.Range("T10").Value = 0.224150728206771

Write-Output "Applied TPM updates to rows 2-10"